# Adds row 5 ("2021年") of data to Sheet1, mirroring the existing yearly rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetRow = 5

# Year label in column A, matching the style used by the existing year cells (A2:A4).
# Copy the previous year cell first so the bold/border/alignment formatting carries
# over, then overwrite the value with the new year label.
$ws.Cells.Item($targetRow - 1, 1).Copy($ws.Cells.Item($targetRow, 1))
$ws.Cells.Item($targetRow, 1).Value = "2021年"

# (column number, value) pairs for columns B..DK. Columns G and AG are intentionally
# left blank (no value) for this row, matching the source data which has no reading
# for those two industries in 2021.
$rowData = @(
    ,(2, 17)
    ,(3, 28.1)
    ,(4, -18.1)
    ,(5, 47.9)
    ,(6, -2.6)
    ,(8, 10.4)
    ,(9, 4.2)
    ,(10, 4)
    ,(11, -10.5)
    ,(12, -5.4)
    ,(13, 14.4)
    ,(14, -51.9)
    ,(15, -6.3)
    ,(16, 10.8)
    ,(17, -11.3)
    ,(18, 29.1)
    ,(19, 36.4)
    ,(20, 2.9)
    ,(21, 5.5)
    ,(22, -3.1)
    ,(23, 17.1)
    ,(24, 13.3)
    ,(25, 15.9)
    ,(26, 33.6)
    ,(27, 15.2)
    ,(28, 9.1)
    ,(29, 12.4)
    ,(30, 3.7)
    ,(31, 9.800000000000001)
    ,(32, -53.9)
    ,(34, -33.9)
    ,(35, -7.2)
    ,(36, -28.1)
    ,(37, 32.3)
    ,(38, -2.6)
    ,(39, -5)
    ,(40, -15.3)
    ,(41, -15.2)
    ,(42, -35.1)
    ,(43, 2.3)
    ,(44, 2.3)
    ,(45, 176.8)
    ,(46, 109.5)
    ,(47, 57.8)
    ,(48, 4.3)
    ,(49, 12.4)
    ,(50, -13.8)
    ,(51, -13.6)
    ,(52, 11.1)
    ,(53, -4.5)
    ,(54, -1.3)
    ,(55, 12.5)
    ,(56, 74.59999999999999)
    ,(57, -0.3)
    ,(58, 7.8)
    ,(59, 10.9)
    ,(60, -20.4)
    ,(61, 2.8)
    ,(62, 13.8)
    ,(63, 16.6)
    ,(64, -6.2)
    ,(65, -3.4)
    ,(66, -9.6)
    ,(67, -6.2)
    ,(68, -0.1)
    ,(69, 38.6)
    ,(70, 11.2)
    ,(71, -4.3)
    ,(72, -7.1)
    ,(73, -5.3)
    ,(74, 1.1)
    ,(75, 4.7)
    ,(76, 26.1)
    ,(77, 15)
    ,(78, 2)
    ,(79, 8.199999999999999)
    ,(80, 2.9)
    ,(81, 13.3)
    ,(82, -18.3)
    ,(83, -7.1)
    ,(84, 15.5)
    ,(85, 16.7)
    ,(86, -3.7)
    ,(87, 9.199999999999999)
    ,(88, -16.7)
    ,(89, 14.9)
    ,(90, 0.9)
    ,(91, -24.6)
    ,(92, 19.4)
    ,(93, 7.9)
    ,(94, 14.4)
    ,(95, 7.3)
    ,(96, -21.3)
    ,(97, 19.4)
    ,(98, 9.800000000000001)
    ,(99, 18.8)
    ,(100, -5.1)
    ,(101, 21.1)
    ,(102, 14.4)
    ,(103, 11.6)
    ,(104, 75.8)
    ,(105, 8.199999999999999)
    ,(106, 5.3)
    ,(107, 16.7)
    ,(108, -6.1)
    ,(109, -13.4)
    ,(110, 16.9)
    ,(111, 31.4)
    ,(112, 7.4)
    ,(113, 1.4)
    ,(114, 16.5)
    ,(115, 26.3)
)

foreach ($pair in $rowData) {
    $col = $pair[0]
    $val = $pair[1]
    $ws.Cells.Item($targetRow, $col).Value = $val
}

Write-Output "Row 5 (2021年) written: $($rowData.Count) data cells."
